$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny floating-point refinement of the existing last row's timestamp
$ws.Range("A87").Value = 44400.76775271643

# New row of data retrieved at commit time (Sat Jul 24 18:29:01 UTC 2021)
$ws.Range("A88").Value = 44401.77015799552
$ws.Range("B88").Value = 80151
$ws.Range("C88").Value = 67617
$ws.Range("D88").Value = 3712
$ws.Range("E88").Value = 2208
$ws.Range("F88").Value = 1585
$ws.Range("G88").Value = 21002
$ws.Range("H88").Value = 1615
$ws.Range("I88").Value = 888
$ws.Range("J88").Value = 203

$ws.Range("A88").NumberFormat = $ws.Range("A87").NumberFormat
